# Fanano.xlsx — "aggiornamento fino a 28 luglio"
# Appends new daily-report rows (date serial, nuovi pos., somma mobile 7gg.,
# somma mobile 7gg. per 100mila abitanti) below the existing data, extending
# the sheet from A1:D301 to A1:D328 (dates 2021-06-29 .. 2021-07-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 302
$lastExistingRow = 301

# New rows: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$data = @(
    @(44376, 0, 0, 0),
    @(44377, 0, 0, 0),
    @(44378, 0, 0, 0),
    @(44379, 1, 1, 33.71544167228591),
    @(44380, 0, 1, 33.71544167228591),
    @(44381, 0, 1, 33.71544167228591),
    @(44382, 0, 1, 33.71544167228591),
    @(44383, 0, 1, 33.71544167228591),
    @(44384, 0, 1, 33.71544167228591),
    @(44385, 0, 1, 33.71544167228591),
    @(44386, 0, 0, 0),
    @(44387, 0, 0, 0),
    @(44388, 0, 0, 0),
    @(44389, 0, 0, 0),
    @(44390, 0, 0, 0),
    @(44391, 0, 0, 0),
    @(44392, 0, 0, 0),
    @(44393, 0, 0, 0),
    @(44394, 0, 0, 0),
    @(44395, 0, 0, 0),
    @(44396, 0, 0, 0),
    @(44397, 0, 0, 0),
    @(44398, 0, 0, 0),
    @(44399, 0, 0, 0),
    @(44400, 0, 0, 0),
    @(44401, 0, 0, 0),
    @(44402, 0, 0, 0)
)

# Column A carries the bold/border/centered "date" cell style (s="2" in the
# sheet XML) applied to every row in the existing table. Copy that format
# (format only, not the value) from the last existing row down across all of
# the new A-column cells before writing the new values/dates into them.
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A$($firstNewRow):A$($firstNewRow + $data.Count - 1)").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstNewRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "Appended $($data.Count) rows ($firstNewRow..$($firstNewRow + $data.Count - 1))"
